# Insert a new weekly price-report row at row 26 (pushing existing rows
# 26-57 down to 27-58) and populate it with the new Ají ("Americana (o)")
# observation for "Región del Maule".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 26:57 down to 27:58, carrying formatting (e.g. the date
# number-format on column D) the same way Excel's own row-insert does.
$ws.Rows.Item(26).Insert()

# Populate the newly-inserted row 26 with the new record.
$ws.Cells.Item(26, 1).Value  = 7
$ws.Cells.Item(26, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(26, 3).Value  = "Ñuble"
$ws.Cells.Item(26, 4).Value  = (Get-Date -Year 2022 -Month 1 -Day 19 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(26, 5).Value  = 16
$ws.Cells.Item(26, 6).Value  = 100112021
$ws.Cells.Item(26, 7).Value  = "Ají"
$ws.Cells.Item(26, 8).Value  = "Americana (o)"
$ws.Cells.Item(26, 9).Value  = "Primera"
$ws.Cells.Item(26, 10).Value = 60
$ws.Cells.Item(26, 11).Value = 16000
$ws.Cells.Item(26, 12).Value = 17000
$ws.Cells.Item(26, 13).Value = 16500
$ws.Cells.Item(26, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(26, 15).Value = "Región del Maule"
$ws.Cells.Item(26, 16).Value = 1100
$ws.Cells.Item(26, 17).Value = 15
$ws.Cells.Item(26, 18).Value = "Hortaliza"
